# Update the "Förändrad" (Changed) date column (C) from 45179 (2023-09-10)
# to 45180 (2023-09-11) for every data row (rows 2 through 121).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 121; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value2 = 45180
    }
}
